$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计" (position 2)
# ------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item(1)
$q4 = $wb.Worksheets.Add($null, $summaryWs)
$q4.Name = "2022-Q4"

$hdr = $q4.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# Populate data rows
$q4.Range("B1").Value = "'基金代码"
$q4.Range("C1").Value = "'基金名称"
$q4.Range("D1").Value = "'基金规模"
$q4.Range("E1").Value = "'股票总仓位"
$q4.Range("F1").Value = "'仓位占比"
$q4.Range("G1").Value = "'持有市值(亿元)"
$q4.Range("H1").Value = "'仓位排名"

$q4.Range("A2").Value = 0
$q4.Range("B2").Value = "'006551"
$q4.Range("C2").Value = "'中庚价值领航混合"
$q4.Range("D2").Value = "'116.63"
$q4.Range("E2").Value = "'93.22"
$q4.Range("F2").Value = "'5.08"
$q4.Range("G2").Value = "'5.9248"
$q4.Range("H2").Value = 6

$q4.Range("A3").Value = 1
$q4.Range("B3").Value = "'011174"
$q4.Range("C3").Value = "'中庚价值品质一年持有期混合"
$q4.Range("D3").Value = "'67.05"
$q4.Range("E3").Value = "'93.59"
$q4.Range("F3").Value = "'5.50"
$q4.Range("G3").Value = "'3.6878"
$q4.Range("H3").Value = 6

$q4.Range("A4").Value = 2
$q4.Range("B4").Value = "'007130"
$q4.Range("C4").Value = "'中庚小盘价值股票"
$q4.Range("D4").Value = "'76.15"
$q4.Range("E4").Value = "'93.50"
$q4.Range("F4").Value = "'4.59"
$q4.Range("G4").Value = "'3.4953"
$q4.Range("H4").Value = 5

$q4.Range("A5").Value = 3
$q4.Range("B5").Value = "'090018"
$q4.Range("C5").Value = "'大成新锐产业混合"
$q4.Range("D5").Value = "'88.75"
$q4.Range("E5").Value = "'93.33"
$q4.Range("F5").Value = "'2.93"
$q4.Range("G5").Value = "'2.6004"
$q4.Range("H5").Value = 10

$q4.Range("A6").Value = 4
$q4.Range("B6").Value = "'001300"
$q4.Range("C6").Value = "'大成睿景灵活配置混合A"
$q4.Range("D6").Value = "'32.95"
$q4.Range("E6").Value = "'92.29"
$q4.Range("F6").Value = "'2.92"
$q4.Range("G6").Value = "'0.9621"
$q4.Range("H6").Value = 10

$q4.Range("A7").Value = 5
$q4.Range("B7").Value = "'013435"
$q4.Range("C7").Value = "'大成景气精选六个月持有混合A"
$q4.Range("D7").Value = "'30.45"
$q4.Range("E7").Value = "'91.16"
$q4.Range("F7").Value = "'2.92"
$q4.Range("G7").Value = "'0.8891"
$q4.Range("H7").Value = 10

$q4.Range("A8").Value = 6
$q4.Range("B8").Value = "'001301"
$q4.Range("C8").Value = "'大成睿景灵活配置混合C"
$q4.Range("D8").Value = "'19.30"
$q4.Range("E8").Value = "'92.29"
$q4.Range("F8").Value = "'2.92"
$q4.Range("G8").Value = "'0.5636"
$q4.Range("H8").Value = 10

$q4.Range("A9").Value = 7
$q4.Range("B9").Value = "'002258"
$q4.Range("C9").Value = "'大成国企改革灵活配置混合"
$q4.Range("D9").Value = "'16.71"
$q4.Range("E9").Value = "'93.37"
$q4.Range("F9").Value = "'2.92"
$q4.Range("G9").Value = "'0.4879"
$q4.Range("H9").Value = 10

$q4.Range("A10").Value = 8
$q4.Range("B10").Value = "'012519"
$q4.Range("C10").Value = "'大成核心趋势混合A"
$q4.Range("D10").Value = "'9.51"
$q4.Range("E10").Value = "'91.16"
$q4.Range("F10").Value = "'2.92"
$q4.Range("G10").Value = "'0.2777"
$q4.Range("H10").Value = 10

$q4.Range("A11").Value = 9
$q4.Range("B11").Value = "'013436"
$q4.Range("C11").Value = "'大成景气精选六个月持有混合C"
$q4.Range("D11").Value = "'5.20"
$q4.Range("E11").Value = "'91.16"
$q4.Range("F11").Value = "'2.92"
$q4.Range("G11").Value = "'0.1518"
$q4.Range("H11").Value = 10

$q4.Range("A12").Value = 10
$q4.Range("B12").Value = "'012520"
$q4.Range("C12").Value = "'大成核心趋势混合C"
$q4.Range("D12").Value = "'2.45"
$q4.Range("E12").Value = "'91.16"
$q4.Range("F12").Value = "'2.92"
$q4.Range("G12").Value = "'0.0715"
$q4.Range("H12").Value = 10

$q4.Range("A13").Value = 11
$q4.Range("B13").Value = "'260117"
$q4.Range("C13").Value = "'景顺长城支柱产业混合"
$q4.Range("D13").Value = "'0.77"
$q4.Range("E13").Value = "'72.88"
$q4.Range("F13").Value = "'3.83"
$q4.Range("G13").Value = "'0.0295"
$q4.Range("H13").Value = 9

$q4.Range("A14").Value = 12
$q4.Range("B14").Value = "'159990"
$q4.Range("C14").Value = "'银华巨潮小盘价值ETF"
$q4.Range("D14").Value = "'0.78"
$q4.Range("E14").Value = "'97.02"
$q4.Range("F14").Value = "'1.17"
$q4.Range("G14").Value = "'0.0091"
$q4.Range("H14").Value = 7

# Style the index column (A2:A14) to match header look
$idxCol = $q4.Range("A2:A14")
$idxCol.Font.Bold = $true
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160
$idxCol.Borders.LineStyle = 1
# ------------------------------------------------------------------
# 2) Update the summary ("总计") sheet: insert a new row for 2022-Q4
#    and shift the rest down by one row.
# ------------------------------------------------------------------
$summaryWs.Rows.Item(2).Insert()

$a2 = $summaryWs.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$summaryWs.Range("B2").Value = "'2022-Q4"
$summaryWs.Range("C2").Value = 13
$summaryWs.Range("D2").Value = 19.15
